$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title text (shared string moved/edited)
$ws.Range("A1").Value = "Analyseresultater fra forsurede og kalkede vassdrag - kvartalsrapport nr. XXX / årsrapport XXXX"

# Update the selection to match the new sqref (A1:AD1), no explicit active cell anchor change needed
$ws.Range("A1:AD1").Select()
